$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Rv1185c", 3, "fadD21 Rv1185c MTV005.21c", "", 16),
    @("Rv1068c", 2, "PE_PGRS20 Rv1068c MTV017.21c", "", 16),
    @("Rv0504c", 2, "Rv0504c MTCY20G9.31c", "", 16),
    @("Rv1357c", 2, "Rv1357c MTCY02B10.21c", "", 16),
    @("Rv1220c", 2, "Rv1220c", "", 16),
    @("Rv1565c", 2, "Rv1565c", "", 16),
    @("Rv3759c", 2, "proX Rv3759c", "", 16),
    @("Rv0481c", 1, "Rv0481c MTCY20G9.07c", "", 16),
    @("Rv2929", 1, "Rv2929 MTCY338.18", "", 16),
    @("Rv2901c", 1, "Rv2901c MTCY274.32c", "", 16),
    @("Rv2360c", 1, "Rv2360c", "", 16),
    @("Rv2425c", 1, "Rv2425c", "", 16),
    @("Rv3093c", 1, "Rv3093c", "", 16)
)

$row = 3
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = ""
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row++
}
